$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new "Surveys" variable definitions (Mnemonic / Description / Group)
# Cells are written in the same order they were originally typed in (Description,
# then Mnemonic, then Group for the first new row; Group is reused afterwards;
# the last row was entered Mnemonic-then-Description) so that the resulting
# shared-strings table matches the authored workbook.

$ws.Cells.Item(30, 2).Value = "ifo: industry and trade, climate"
$ws.Cells.Item(30, 1).Value = "ifoIndTradeClimate"
$ws.Cells.Item(30, 3).Value = "Surveys"

$ws.Cells.Item(31, 2).Value = "ifo: industry and trade, current situation"
$ws.Cells.Item(31, 1).Value = "ifoIndTradeCurrent"
$ws.Cells.Item(31, 3).Value = "Surveys"

$ws.Cells.Item(32, 2).Value = "ifo: industry and trade, expectations"
$ws.Cells.Item(32, 1).Value = "ifoIndTradeExp"
$ws.Cells.Item(32, 3).Value = "Surveys"

$ws.Cells.Item(33, 2).Value = "GfK: business cycle expectations"
$ws.Cells.Item(33, 1).Value = "GfKBCE"
$ws.Cells.Item(33, 3).Value = "Surveys"

$ws.Cells.Item(34, 2).Value = "GfK: income expectations"
$ws.Cells.Item(34, 1).Value = "GfKIE"
$ws.Cells.Item(34, 3).Value = "Surveys"

$ws.Cells.Item(35, 2).Value = "GfK: willigness-to-buy"
$ws.Cells.Item(35, 1).Value = "GfKWtB"
$ws.Cells.Item(35, 3).Value = "Surveys"

$ws.Cells.Item(36, 2).Value = "GfK: consumer climate indicator"
$ws.Cells.Item(36, 1).Value = "GfKCCI"
$ws.Cells.Item(36, 3).Value = "Surveys"

$ws.Cells.Item(37, 1).Value = "ESI"
$ws.Cells.Item(37, 2).Value = "Economics Sentiment Indicator"
$ws.Cells.Item(37, 3).Value = "Surveys"

# Update the saved selection / scroll position to match the author's view
$ws.Range("B41").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
